$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.823505401611328
$ws.Range("B1").Value = 2.02829384803772
$ws.Range("C1").Value = 2.213552474975586
$ws.Range("D1").Value = 3.24401330947876
$ws.Range("E1").Value = 1.884402513504028
